$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("vaccines")
$ws.Activate()

# Rows 4 (AstraZeneca) and 6 (Curevac): status_en/status_de/status_fr were
# changed from "approval process aborted" / "Zulassungsverfahren
# abgebrochen" / "procédure d'approbation interrompue" to
# "not approved" / "nicht zugelassen" / "non autorisé".
$ws.Range("AE4").Value = "not approved"
$ws.Range("AF4").Value = "nicht zugelassen"
$ws.Range("AG4").Value = "non autorisé"

$ws.Range("AE6").Value = "not approved"
$ws.Range("AF6").Value = "nicht zugelassen"
$ws.Range("AG6").Value = "non autorisé"

# Widen the now more-verbose status_en / status_de columns.
$ws.Columns("AE").ColumnWidth = 21.17
$ws.Columns("AF").ColumnWidth = 26

# Move the current selection on the frozen pane to X13.
$ws.Range("X13").Select()
